# 1st commit on 3rd July 2021
# Rename sheet "8" -> "7", add a new sheet "10" with a second trade row,
# and extend sheet "7" with a merged second trade (rows 2-3).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "7" (was "8")
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "7"

# Update existing row 2 (first trade row) with the new trade data.
$ws1.Range("A2").Value = 6
$ws1.Range("B2").Value = "VNDUSD"
$ws1.Range("C2").Value = "14:57:30.000793"
$ws1.Range("D2").Value = "Sell"
$ws1.Range("E2").Value = "link 1"
$ws1.Range("F2").Value = "link 2"
$ws1.Range("G2").Value = "link 3"
$ws1.Range("H2").Value = "link 4"
$ws1.Range("I2").Value = "link 5"
$ws1.Range("J2").Value = 3
$ws1.Range("K2").Value = "this is my comment"
$ws1.Range("L2").Value = 606

# Make A2 store as literal text ("6"), matching the source workbook,
# instead of Excel's default numeric auto-detection.
$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "6"

# Copy the formatting of row 2 down into the new row 3 before filling it in.
$ws1.Range("A2:M2").Copy($ws1.Range("A3:M3"))

# Row 3 (second trade, sharing DAY/ID + SUM with row 2 via the merge below).
$ws1.Range("B3").Value = "PAIR2"
$ws1.Range("C3").Value = "14:58:01.136710"
$ws1.Range("D3").Value = "Sell"
$ws1.Range("E3").Value = "link 1"
$ws1.Range("F3").Value = "link 2"
$ws1.Range("G3").Value = "link 3"
$ws1.Range("H3").Value = "link 4"
$ws1.Range("I3").Value = "link 5"
$ws1.Range("J3").Value = 3
$ws1.Range("K3").Value = "this is my comment"
$ws1.Range("L3").Value = 806

# SUM column becomes a formula totalling the two trade rows.
$ws1.Range("M2").Formula = "=SUM(J2:J3)"
$ws1.Range("M3").Value = ""

# Merge the DAY/ID cell and the SUM cell across the two trade rows.
$ws1.Range("A2:A3").Merge()
$ws1.Range("M2:M3").Merge()

# ---------------------------------------------------------------------
# Sheet "10" (brand new sheet, placed right after "7")
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "10"

# Bring over the header row (and its formatting) from sheet "7".
$ws1.Range("A1:M1").Copy($ws2.Range("A1:M1"))

# Approximate the wide "link" columns / comments column widths from sheet "7".
$ws2.Range("E1:I1").ColumnWidth = 34.8
$ws2.Range("K1").ColumnWidth = 29.8

# Give row 2 the same look (borders/alignment) as the header before filling it in.
$ws1.Range("A1:M1").Copy($ws2.Range("A2:M2"))

$ws2.Range("A2").Value = 6
$ws2.Range("B2").Value = "VNDUSD"
$ws2.Range("C2").Value = "14:57:03.901108"
$ws2.Range("D2").Value = "Sell"
$ws2.Range("E2").Value = "link 1"
$ws2.Range("F2").Value = "link 2"
$ws2.Range("G2").Value = "link 3"
$ws2.Range("H2").Value = "link 4"
$ws2.Range("I2").Value = "link 5"
$ws2.Range("J2").Value = 3
$ws2.Range("K2").Value = "this is my comment"
$ws2.Range("L2").Value = 406
$ws2.Range("M2").Value = 3

$ws1.Select()
